# Refresh the cryptocurrency price/volume snapshot (coinranking.com feed).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCellValue {
    param($Sheet, $Row, $Col, $Text)
    $cell = $Sheet.Cells.Item($Row, $Col)
    # Force text storage so numeric-looking strings (e.g. "0.9985") are not
    # auto-coerced to real numbers by Excel's input parser, then restore the
    # default "Normal" style so no stray number-format is left on the cell.
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.NumberFormat = "General"
    $cell.Style = "Normal"
}

# Price (column D) / Volume(1h) (column E) values refreshed this run.
$updates = @(
    @{Row=2; Col=4; Val='27.954.40'},
    @{Row=2; Col=5; Val='  -0.26%  '},
    @{Row=3; Col=4; Val='1.910.41'},
    @{Row=3; Col=5; Val='  -0.07%  '},
    @{Row=4; Col=4; Val='0.9985'},
    @{Row=4; Col=5; Val='  -0.70%  '},
    @{Row=5; Col=4; Val='313.65'},
    @{Row=5; Col=5; Val='  -0.67%  '},
    @{Row=6; Col=4; Val='0.9985'},
    @{Row=7; Col=4; Val='0.5007'},
    @{Row=7; Col=5; Val='  +3.93%  '},
    @{Row=8; Col=4; Val='0.3827'},
    @{Row=8; Col=5; Val='  +0.38%  '},
    @{Row=9; Col=4; Val='0.07329'},
    @{Row=9; Col=5; Val='  -0.42%  '},
    @{Row=10; Col=4; Val='0.9138'},
    @{Row=10; Col=5; Val='  -2.13%  '},
    @{Row=11; Col=4; Val='21.18'},
    @{Row=11; Col=5; Val='  +1.79%  '},
    @{Row=12; Col=4; Val='0.07683'},
    @{Row=12; Col=5; Val='  -1.33%  '},
    @{Row=13; Col=4; Val='1.905.74'},
    @{Row=13; Col=5; Val='  -0.37%  '},
    @{Row=14; Col=4; Val='5.518'},
    @{Row=14; Col=5; Val='  +0.24%  '},
    @{Row=15; Col=4; Val='92.79'},
    @{Row=15; Col=5; Val='  +0.90%  '},
    @{Row=16; Col=4; Val='0.9986'},
    @{Row=16; Col=5; Val='  -0.81%  '},
    @{Row=17; Col=4; Val='0.000008755'},
    @{Row=17; Col=5; Val='  -1.29%  '},
    @{Row=18; Col=4; Val='0.9977'},
    @{Row=18; Col=5; Val='  -0.68%  '},
    @{Row=19; Col=4; Val='27.983.42'},
    @{Row=19; Col=5; Val='  -0.29%  '},
    @{Row=20; Col=4; Val='14.68'},
    @{Row=20; Col=5; Val='  -0.65%  '},
    @{Row=21; Col=4; Val='5.187'},
    @{Row=21; Col=5; Val='  +0.31%  '},
    @{Row=22; Col=4; Val='10.86'},
    @{Row=22; Col=5; Val='  -0.53%  '},
    @{Row=23; Col=4; Val='6.608'},
    @{Row=23; Col=5; Val='  -0.56%  '},
    @{Row=24; Col=5; Val='  -2.02%  '},
    @{Row=25; Col=4; Val='1.845'},
    @{Row=25; Col=5; Val='  -3.63%  '},
    @{Row=26; Col=4; Val='2.218'},
    @{Row=26; Col=5; Val='  +4.23%  '},
    @{Row=27; Col=5; Val='  -0.32%  '},
    @{Row=28; Col=4; Val='115.65'},
    @{Row=28; Col=5; Val='  -0.98%  '},
    @{Row=29; Col=4; Val='4.920'},
    @{Row=29; Col=5; Val='  -0.94%  '},
    @{Row=30; Col=4; Val='0.09031'},
    @{Row=30; Col=5; Val='  +0.81%  '},
    @{Row=31; Col=4; Val='3.208'},
    @{Row=31; Col=5; Val='  -2.94%  '},
    @{Row=32; Col=4; Val='4.863'},
    @{Row=32; Col=5; Val='  +4.02%  '},
    @{Row=33; Col=4; Val='1.240'},
    @{Row=33; Col=5; Val='  -1.69%  '},
    @{Row=34; Col=4; Val='0.7777'},
    @{Row=34; Col=5; Val='  +0.00%  '},
    @{Row=35; Col=4; Val='0.02091'},
    @{Row=35; Col=5; Val='  +1.78%  '},
    @{Row=36; Col=4; Val='2.588'},
    @{Row=36; Col=5; Val='  -0.93%  '},
    @{Row=37; Col=4; Val='3.074'},
    @{Row=37; Col=5; Val='  +2.71%  '},
    @{Row=38; Col=4; Val='1.093'},
    @{Row=38; Col=5; Val='  -1.65%  '},
    @{Row=39; Col=4; Val='0.5569'},
    @{Row=39; Col=5; Val='  +1.03%  '},
    @{Row=40; Col=4; Val='0.05292'},
    @{Row=40; Col=5; Val='  -0.52%  '},
    @{Row=41; Col=4; Val='6.896'},
    @{Row=41; Col=5; Val='  -1.90%  '},
    @{Row=42; Col=4; Val='8.542'},
    @{Row=42; Col=5; Val='  +0.73%  '},
    @{Row=43; Col=4; Val='112.98'},
    @{Row=43; Col=5; Val='  +3.81%  '},
    @{Row=44; Col=4; Val='0.1524'},
    @{Row=44; Col=5; Val='  -0.31%  '},
    @{Row=47; Col=4; Val='0.9980'},
    @{Row=47; Col=5; Val='  -0.72%  '},
    @{Row=48; Col=4; Val='1.642'},
    @{Row=48; Col=5; Val='  -0.32%  '},
    @{Row=49; Col=4; Val='67.68'},
    @{Row=49; Col=5; Val='  -0.96%  '},
    @{Row=50; Col=4; Val='0.06049'},
    @{Row=50; Col=5; Val='  -0.52%  '},
    @{Row=51; Col=4; Val='0.9087'},
    @{Row=51; Col=5; Val='  +0.96%  '}
)

foreach ($u in $updates) {
    Set-TextCellValue $ws $u.Row $u.Col $u.Val
}

# Rows 45 and 46 swapped rank order this run: Decentraland moved above
# EnergySwap, each carrying its own refreshed price/volume figures.
Set-TextCellValue $ws 45 2 'Decentraland'
Set-TextCellValue $ws 45 3 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextCellValue $ws 45 4 '0.4843'
Set-TextCellValue $ws 45 5 '  +0.36%  '

Set-TextCellValue $ws 46 2 'EnergySwap'
Set-TextCellValue $ws 46 3 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCellValue $ws 46 4 '10.60'
Set-TextCellValue $ws 46 5 '  -0.56%  '
